$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Default New Item Name" - update the template/default new-item row (row 2)
# Nombre / Descripcion (shared text "F030-PANTALONETA" -> "SQ009-DELINEADOR")
$ws.Range("A2").Value = "SQ009-DELINEADOR"
$ws.Range("P2").Value = "SQ009-DELINEADOR"

# Codigo Interno / Cod barras (shared text "F030" -> "SQ009")
$ws.Range("B2").Value = "SQ009"
$ws.Range("T2").Value = "SQ009"

# "Val" - Precio Unitario Venta changed from 70 to 12
$ws.Range("G2").Value = 12
